$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the substantive input values (B13, D13)
$ws.Range("B13").Value = 1.1200000000000001
$ws.Range("D13").Value = 0.18

# Change B13 number format from scientific (0.00E+00) to General
$ws.Range("B13").NumberFormat = "General"

# Update the selected cell to B17
$ws.Range("B17").Select()
